$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D shifts to F, etc.)
$ws.Columns("D:E").Insert()

# Copy number formatting from column F (the shifted former column D) into the new D:E columns
# so the new columns inherit the correct style (date format for header rows, number format for data rows).
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Write the refreshed financial figures (new quarters in D:E, corrected historical figures in F:M)
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(7, 6).Value = 43281
$ws.Cells.Item(7, 7).Value = 43190
$ws.Cells.Item(7, 8).Value = 43100
$ws.Cells.Item(7, 9).Value = 43008
$ws.Cells.Item(7, 10).Value = 42916
$ws.Cells.Item(7, 11).Value = 42825
$ws.Cells.Item(7, 12).Value = 42735
$ws.Cells.Item(7, 13).Value = 42643
$ws.Cells.Item(8, 4).Value = 98700
$ws.Cells.Item(8, 5).Value = 124900
$ws.Cells.Item(8, 6).Value = 99100
$ws.Cells.Item(8, 7).Value = 75000
$ws.Cells.Item(8, 8).Value = 93500
$ws.Cells.Item(8, 9).Value = 124600
$ws.Cells.Item(8, 10).Value = 102100
$ws.Cells.Item(8, 11).Value = 69000
$ws.Cells.Item(8, 12).Value = 79300
$ws.Cells.Item(8, 13).Value = 112300
$ws.Cells.Item(9, 4).Value = 47600
$ws.Cells.Item(9, 5).Value = 60100
$ws.Cells.Item(9, 6).Value = 44500
$ws.Cells.Item(9, 7).Value = 34900
$ws.Cells.Item(9, 8).Value = 44900
$ws.Cells.Item(9, 9).Value = 57900
$ws.Cells.Item(9, 10).Value = 43300
$ws.Cells.Item(9, 11).Value = 29800
$ws.Cells.Item(9, 12).Value = 36400
$ws.Cells.Item(9, 13).Value = 49200
$ws.Cells.Item(10, 4).Value = 51100
$ws.Cells.Item(10, 5).Value = 64800
$ws.Cells.Item(10, 6).Value = 54600
$ws.Cells.Item(10, 7).Value = 40100
$ws.Cells.Item(10, 8).Value = 48600
$ws.Cells.Item(10, 9).Value = 66700
$ws.Cells.Item(10, 10).Value = 58800
$ws.Cells.Item(10, 11).Value = 39200
$ws.Cells.Item(10, 12).Value = 42900
$ws.Cells.Item(10, 13).Value = 63100
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(12, 6).Value = "NA"
$ws.Cells.Item(12, 7).Value = "NA"
$ws.Cells.Item(12, 8).Value = "NA"
$ws.Cells.Item(12, 9).Value = "NA"
$ws.Cells.Item(12, 10).Value = "NA"
$ws.Cells.Item(12, 11).Value = "NA"
$ws.Cells.Item(12, 12).Value = "NA"
$ws.Cells.Item(12, 13).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = 0
$ws.Cells.Item(14, 4).Value = 3600
$ws.Cells.Item(14, 5).Value = 8400
$ws.Cells.Item(14, 6).Value = 2700
$ws.Cells.Item(14, 7).Value = 3800
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = "NA"
$ws.Cells.Item(14, 12).Value = "NA"
$ws.Cells.Item(14, 13).Value = "NA"
$ws.Cells.Item(15, 4).Value = 14300
$ws.Cells.Item(15, 5).Value = 14300
$ws.Cells.Item(15, 6).Value = 14200
$ws.Cells.Item(15, 7).Value = 14200
$ws.Cells.Item(15, 8).Value = 12600
$ws.Cells.Item(15, 9).Value = 12600
$ws.Cells.Item(15, 10).Value = 12600
$ws.Cells.Item(15, 11).Value = 12600
$ws.Cells.Item(15, 12).Value = 11600
$ws.Cells.Item(15, 13).Value = 11500
$ws.Cells.Item(17, 4).Value = 81200
$ws.Cells.Item(17, 5).Value = 99000
$ws.Cells.Item(17, 6).Value = 76300
$ws.Cells.Item(17, 7).Value = 67700
$ws.Cells.Item(17, 8).Value = 74400
$ws.Cells.Item(17, 9).Value = 85400
$ws.Cells.Item(17, 10).Value = 70300
$ws.Cells.Item(17, 11).Value = 58200
$ws.Cells.Item(17, 12).Value = 62900
$ws.Cells.Item(17, 13).Value = 76000
$ws.Cells.Item(18, 4).Value = 17500
$ws.Cells.Item(18, 5).Value = 25900
$ws.Cells.Item(18, 6).Value = 22800
$ws.Cells.Item(18, 7).Value = 7300
$ws.Cells.Item(18, 8).Value = 19100
$ws.Cells.Item(18, 9).Value = 39200
$ws.Cells.Item(18, 10).Value = 31800
$ws.Cells.Item(18, 11).Value = 10800
$ws.Cells.Item(18, 12).Value = 16400
$ws.Cells.Item(18, 13).Value = 36300
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(20, 5).Value = 100
$ws.Cells.Item(20, 6).Value = 200
$ws.Cells.Item(20, 7).Value = -500
$ws.Cells.Item(20, 8).Value = 12200
$ws.Cells.Item(20, 9).Value = -600
$ws.Cells.Item(20, 10).Value = 6500
$ws.Cells.Item(20, 11).Value = 500
$ws.Cells.Item(20, 12).Value = 10800
$ws.Cells.Item(20, 13).Value = 600
$ws.Cells.Item(21, 4).Value = 31700
$ws.Cells.Item(21, 5).Value = 40200
$ws.Cells.Item(21, 6).Value = 37300
$ws.Cells.Item(21, 7).Value = 20900
$ws.Cells.Item(21, 8).Value = 44000
$ws.Cells.Item(21, 9).Value = 51300
$ws.Cells.Item(21, 10).Value = 51000
$ws.Cells.Item(21, 11).Value = 23900
$ws.Cells.Item(21, 12).Value = 38800
$ws.Cells.Item(21, 13).Value = 48400
$ws.Cells.Item(22, 4).Value = 6100
$ws.Cells.Item(22, 5).Value = 6100
$ws.Cells.Item(22, 6).Value = 6100
$ws.Cells.Item(22, 7).Value = 6100
$ws.Cells.Item(22, 8).Value = 5700
$ws.Cells.Item(22, 9).Value = 5600
$ws.Cells.Item(22, 10).Value = 5800
$ws.Cells.Item(22, 11).Value = 6100
$ws.Cells.Item(22, 12).Value = 5500
$ws.Cells.Item(22, 13).Value = 5400
$ws.Cells.Item(23, 4).Value = 11300
$ws.Cells.Item(23, 5).Value = 19900
$ws.Cells.Item(23, 6).Value = 16900
$ws.Cells.Item(23, 7).Value = 700
$ws.Cells.Item(23, 8).Value = 25600
$ws.Cells.Item(23, 9).Value = 33100
$ws.Cells.Item(23, 10).Value = 32500
$ws.Cells.Item(23, 11).Value = 5200
$ws.Cells.Item(23, 12).Value = 21700
$ws.Cells.Item(23, 13).Value = 31500
$ws.Cells.Item(24, 4).Value = 2500
$ws.Cells.Item(24, 5).Value = 4100
$ws.Cells.Item(24, 6).Value = 4100
$ws.Cells.Item(24, 7).Value = -600
$ws.Cells.Item(24, 8).Value = 10700
$ws.Cells.Item(24, 9).Value = 13500
$ws.Cells.Item(24, 10).Value = 12000
$ws.Cells.Item(24, 11).Value = 1600
$ws.Cells.Item(24, 12).Value = 8000
$ws.Cells.Item(24, 13).Value = 12500
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = 0
$ws.Cells.Item(26, 4).Value = 8800
$ws.Cells.Item(26, 5).Value = 15800
$ws.Cells.Item(26, 6).Value = 12900
$ws.Cells.Item(26, 7).Value = 1300
$ws.Cells.Item(26, 8).Value = 14900
$ws.Cells.Item(26, 9).Value = 19500
$ws.Cells.Item(26, 10).Value = 20600
$ws.Cells.Item(26, 11).Value = 3700
$ws.Cells.Item(26, 12).Value = 13700
$ws.Cells.Item(26, 13).Value = 19000
$ws.Cells.Item(27, 4).Value = 8800
$ws.Cells.Item(27, 5).Value = 15800
$ws.Cells.Item(27, 6).Value = 12900
$ws.Cells.Item(27, 7).Value = 1300
$ws.Cells.Item(27, 8).Value = 14900
$ws.Cells.Item(27, 9).Value = 19500
$ws.Cells.Item(27, 10).Value = 18700
$ws.Cells.Item(27, 11).Value = 3700
$ws.Cells.Item(27, 12).Value = 13700
$ws.Cells.Item(27, 13).Value = 19000
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 11).Value = 0
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 13).Value = 0
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = "NA"
$ws.Cells.Item(29, 6).Value = "NA"
$ws.Cells.Item(29, 7).Value = "NA"
$ws.Cells.Item(29, 8).Value = 2400
$ws.Cells.Item(29, 9).Value = "NA"
$ws.Cells.Item(29, 10).Value = "NA"
$ws.Cells.Item(29, 11).Value = "NA"
$ws.Cells.Item(29, 12).Value = "NA"
$ws.Cells.Item(29, 13).Value = "NA"
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = 0
$ws.Cells.Item(32, 4).Value = 0
$ws.Cells.Item(32, 5).Value = -100
$ws.Cells.Item(32, 6).Value = -200
$ws.Cells.Item(32, 7).Value = 500
$ws.Cells.Item(32, 8).Value = -12200
$ws.Cells.Item(32, 9).Value = 600
$ws.Cells.Item(32, 10).Value = -6500
$ws.Cells.Item(32, 11).Value = -500
$ws.Cells.Item(32, 12).Value = -10800
$ws.Cells.Item(32, 13).Value = -600
$ws.Cells.Item(33, 4).Value = 8800
$ws.Cells.Item(33, 5).Value = 15800
$ws.Cells.Item(33, 6).Value = 12900
$ws.Cells.Item(33, 7).Value = 1300
$ws.Cells.Item(33, 8).Value = 17300
$ws.Cells.Item(33, 9).Value = 19500
$ws.Cells.Item(33, 10).Value = 18700
$ws.Cells.Item(33, 11).Value = 3700
$ws.Cells.Item(33, 12).Value = 13700
$ws.Cells.Item(33, 13).Value = 19000
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = 0
$ws.Cells.Item(35, 4).Value = 8800
$ws.Cells.Item(35, 5).Value = 15800
$ws.Cells.Item(35, 6).Value = 12900
$ws.Cells.Item(35, 7).Value = 1300
$ws.Cells.Item(35, 8).Value = 17300
$ws.Cells.Item(35, 9).Value = 19500
$ws.Cells.Item(35, 10).Value = 18700
$ws.Cells.Item(35, 11).Value = 3700
$ws.Cells.Item(35, 12).Value = 13700
$ws.Cells.Item(35, 13).Value = 19000
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(38, 6).Value = 43281
$ws.Cells.Item(38, 7).Value = 43190
$ws.Cells.Item(38, 8).Value = 43100
$ws.Cells.Item(38, 9).Value = 43008
$ws.Cells.Item(38, 10).Value = 42916
$ws.Cells.Item(38, 11).Value = 42825
$ws.Cells.Item(38, 12).Value = 42735
$ws.Cells.Item(38, 13).Value = 42643
$ws.Cells.Item(41, 4).Value = 420700
$ws.Cells.Item(41, 5).Value = 13300
$ws.Cells.Item(41, 6).Value = 8900
$ws.Cells.Item(41, 7).Value = 7000
$ws.Cells.Item(41, 8).Value = 7800
$ws.Cells.Item(41, 9).Value = 7600
$ws.Cells.Item(41, 10).Value = 9200
$ws.Cells.Item(41, 11).Value = 7100
$ws.Cells.Item(41, 12).Value = 6300
$ws.Cells.Item(41, 13).Value = 17700
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(42, 7).Value = 0
$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 13).Value = 0
$ws.Cells.Item(43, 4).Value = 22200
$ws.Cells.Item(43, 5).Value = 23400
$ws.Cells.Item(43, 6).Value = 25000
$ws.Cells.Item(43, 7).Value = 27000
$ws.Cells.Item(43, 8).Value = 26400
$ws.Cells.Item(43, 9).Value = 24900
$ws.Cells.Item(43, 10).Value = 20700
$ws.Cells.Item(43, 11).Value = 28500
$ws.Cells.Item(43, 12).Value = 29500
$ws.Cells.Item(43, 13).Value = 42400
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(44, 5).Value = 0
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 0
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 13).Value = 0
$ws.Cells.Item(45, 4).Value = 59800
$ws.Cells.Item(45, 5).Value = 43600
$ws.Cells.Item(45, 6).Value = 39300
$ws.Cells.Item(45, 7).Value = 27100
$ws.Cells.Item(45, 8).Value = 32700
$ws.Cells.Item(45, 9).Value = 54600
$ws.Cells.Item(45, 10).Value = 43000
$ws.Cells.Item(45, 11).Value = 38200
$ws.Cells.Item(45, 12).Value = 63700
$ws.Cells.Item(45, 13).Value = 73300
$ws.Cells.Item(46, 4).Value = 502700
$ws.Cells.Item(46, 5).Value = 80400
$ws.Cells.Item(46, 6).Value = 73200
$ws.Cells.Item(46, 7).Value = 61100
$ws.Cells.Item(46, 8).Value = 66900
$ws.Cells.Item(46, 9).Value = 87100
$ws.Cells.Item(46, 10).Value = 72900
$ws.Cells.Item(46, 11).Value = 73700
$ws.Cells.Item(46, 12).Value = 99600
$ws.Cells.Item(46, 13).Value = 133300
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(47, 6).Value = 3200
$ws.Cells.Item(47, 7).Value = 3100
$ws.Cells.Item(47, 8).Value = 4500
$ws.Cells.Item(47, 9).Value = 3800
$ws.Cells.Item(47, 10).Value = 3700
$ws.Cells.Item(47, 11).Value = 3600
$ws.Cells.Item(47, 12).Value = 3400
$ws.Cells.Item(47, 13).Value = 3200
$ws.Cells.Item(48, 4).Value = 1365800
$ws.Cells.Item(48, 5).Value = 1331800
$ws.Cells.Item(48, 6).Value = 1307000
$ws.Cells.Item(48, 7).Value = 1284200
$ws.Cells.Item(48, 8).Value = 1278000
$ws.Cells.Item(48, 9).Value = 1246800
$ws.Cells.Item(48, 10).Value = 1214700
$ws.Cells.Item(48, 11).Value = 1186200
$ws.Cells.Item(48, 12).Value = 1172800
$ws.Cells.Item(48, 13).Value = 1152400
$ws.Cells.Item(49, 4).Value = 7100
$ws.Cells.Item(49, 5).Value = 15700
$ws.Cells.Item(49, 6).Value = 15700
$ws.Cells.Item(49, 7).Value = 15600
$ws.Cells.Item(49, 8).Value = 6300
$ws.Cells.Item(49, 9).Value = 25200
$ws.Cells.Item(49, 10).Value = 25200
$ws.Cells.Item(49, 11).Value = 25100
$ws.Cells.Item(49, 12).Value = 24000
$ws.Cells.Item(49, 13).Value = 23800
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 7).Value = 0
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 13).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = 0
$ws.Cells.Item(52, 4).Value = 80800
$ws.Cells.Item(52, 5).Value = 108900
$ws.Cells.Item(52, 6).Value = 101100
$ws.Cells.Item(52, 7).Value = 100000
$ws.Cells.Item(52, 8).Value = 102300
$ws.Cells.Item(52, 9).Value = 148700
$ws.Cells.Item(52, 10).Value = 154800
$ws.Cells.Item(52, 11).Value = 147200
$ws.Cells.Item(52, 12).Value = 143600
$ws.Cells.Item(52, 13).Value = 141600
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 11).Value = 0
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 13).Value = 0
$ws.Cells.Item(54, 4).Value = 1956400
$ws.Cells.Item(54, 5).Value = 1536900
$ws.Cells.Item(54, 6).Value = 1500100
$ws.Cells.Item(54, 7).Value = 1463900
$ws.Cells.Item(54, 8).Value = 1458000
$ws.Cells.Item(54, 9).Value = 1511600
$ws.Cells.Item(54, 10).Value = 1471200
$ws.Cells.Item(54, 11).Value = 1435800
$ws.Cells.Item(54, 12).Value = 1443400
$ws.Cells.Item(54, 13).Value = 1454300
$ws.Cells.Item(57, 4).Value = 24900
$ws.Cells.Item(57, 5).Value = 27000
$ws.Cells.Item(57, 6).Value = 26200
$ws.Cells.Item(57, 7).Value = 22500
$ws.Cells.Item(57, 8).Value = 23000
$ws.Cells.Item(57, 9).Value = 30700
$ws.Cells.Item(57, 10).Value = 27100
$ws.Cells.Item(57, 11).Value = 20700
$ws.Cells.Item(57, 12).Value = 18700
$ws.Cells.Item(57, 13).Value = 21700
$ws.Cells.Item(58, 4).Value = 100000
$ws.Cells.Item(58, 5).Value = 76000
$ws.Cells.Item(58, 6).Value = 59000
$ws.Cells.Item(58, 7).Value = 39000
$ws.Cells.Item(58, 8).Value = 25000
$ws.Cells.Item(58, 9).Value = 13000
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 11).Value = 100
$ws.Cells.Item(58, 12).Value = 14300
$ws.Cells.Item(58, 13).Value = 75500
$ws.Cells.Item(59, 4).Value = 39000
$ws.Cells.Item(59, 5).Value = 48700
$ws.Cells.Item(59, 6).Value = 39400
$ws.Cells.Item(59, 7).Value = 33900
$ws.Cells.Item(59, 8).Value = 37100
$ws.Cells.Item(59, 9).Value = 49800
$ws.Cells.Item(59, 10).Value = 43300
$ws.Cells.Item(59, 11).Value = 29900
$ws.Cells.Item(59, 12).Value = 30500
$ws.Cells.Item(59, 13).Value = 54000
$ws.Cells.Item(60, 4).Value = 164000
$ws.Cells.Item(60, 5).Value = 151700
$ws.Cells.Item(60, 6).Value = 124600
$ws.Cells.Item(60, 7).Value = 95400
$ws.Cells.Item(60, 8).Value = 85100
$ws.Cells.Item(60, 9).Value = 93400
$ws.Cells.Item(60, 10).Value = 70400
$ws.Cells.Item(60, 11).Value = 50800
$ws.Cells.Item(60, 12).Value = 63600
$ws.Cells.Item(60, 13).Value = 151100
$ws.Cells.Item(61, 4).Value = 431400
$ws.Cells.Item(61, 5).Value = 431300
$ws.Cells.Item(61, 6).Value = 431300
$ws.Cells.Item(61, 7).Value = 431200
$ws.Cells.Item(61, 8).Value = 431100
$ws.Cells.Item(61, 9).Value = 431000
$ws.Cells.Item(61, 10).Value = 430900
$ws.Cells.Item(61, 11).Value = 433400
$ws.Cells.Item(61, 12).Value = 433300
$ws.Cells.Item(61, 13).Value = 364200
$ws.Cells.Item(62, 4).Value = 471700
$ws.Cells.Item(62, 5).Value = 478900
$ws.Cells.Item(62, 6).Value = 479000
$ws.Cells.Item(62, 7).Value = 479500
$ws.Cells.Item(62, 8).Value = 478600
$ws.Cells.Item(62, 9).Value = 534700
$ws.Cells.Item(62, 10).Value = 533800
$ws.Cells.Item(62, 11).Value = 530200
$ws.Cells.Item(62, 12).Value = 524800
$ws.Cells.Item(62, 13).Value = 527400
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 13).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = 0
$ws.Cells.Item(66, 4).Value = 1067100
$ws.Cells.Item(66, 5).Value = 1061900
$ws.Cells.Item(66, 6).Value = 1034800
$ws.Cells.Item(66, 7).Value = 1006000
$ws.Cells.Item(66, 8).Value = 994800
$ws.Cells.Item(66, 9).Value = 1059100
$ws.Cells.Item(66, 10).Value = 1035100
$ws.Cells.Item(66, 11).Value = 1014400
$ws.Cells.Item(66, 12).Value = 1021700
$ws.Cells.Item(66, 13).Value = 1042700
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(69, 6).Value = 0
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 13).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).Value = 0
$ws.Cells.Item(72, 4).Value = 393900
$ws.Cells.Item(72, 5).Value = 390900
$ws.Cells.Item(72, 6).Value = 380900
$ws.Cells.Item(72, 7).Value = 373800
$ws.Cells.Item(72, 8).Value = 376100
$ws.Cells.Item(72, 9).Value = 366800
$ws.Cells.Item(72, 10).Value = 351800
$ws.Cells.Item(72, 11).Value = 337600
$ws.Cells.Item(72, 12).Value = 338400
$ws.Cells.Item(72, 13).Value = 328800
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 0
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 13).Value = 0
$ws.Cells.Item(76, 4).Value = 889300
$ws.Cells.Item(76, 5).Value = 475000
$ws.Cells.Item(76, 6).Value = 465300
$ws.Cells.Item(76, 7).Value = 457800
$ws.Cells.Item(76, 8).Value = 463200
$ws.Cells.Item(76, 9).Value = 452500
$ws.Cells.Item(76, 10).Value = 436100
$ws.Cells.Item(76, 11).Value = 421400
$ws.Cells.Item(76, 12).Value = 421600
$ws.Cells.Item(76, 13).Value = 411600
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(80, 6).Value = 43281
$ws.Cells.Item(80, 7).Value = 43190
$ws.Cells.Item(80, 8).Value = 43100
$ws.Cells.Item(80, 9).Value = 43008
$ws.Cells.Item(80, 10).Value = 42916
$ws.Cells.Item(80, 11).Value = 42825
$ws.Cells.Item(80, 12).Value = 42735
$ws.Cells.Item(80, 13).Value = 42643
$ws.Cells.Item(81, 4).Value = 8800
$ws.Cells.Item(81, 5).Value = 15800
$ws.Cells.Item(81, 6).Value = 12900
$ws.Cells.Item(81, 7).Value = 1300
$ws.Cells.Item(81, 8).Value = 17300
$ws.Cells.Item(81, 9).Value = 19500
$ws.Cells.Item(81, 10).Value = 18700
$ws.Cells.Item(81, 11).Value = 3700
$ws.Cells.Item(81, 12).Value = 13700
$ws.Cells.Item(81, 13).Value = 19000
$ws.Cells.Item(83, 4).Value = 14300
$ws.Cells.Item(83, 5).Value = 14300
$ws.Cells.Item(83, 6).Value = 14200
$ws.Cells.Item(83, 7).Value = 14200
$ws.Cells.Item(83, 8).Value = 12600
$ws.Cells.Item(83, 9).Value = 12600
$ws.Cells.Item(83, 10).Value = 12600
$ws.Cells.Item(83, 11).Value = 12600
$ws.Cells.Item(83, 12).Value = 11600
$ws.Cells.Item(83, 13).Value = 11500
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(84, 6).Value = 0
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 13).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 13).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 13).Value = 0
$ws.Cells.Item(89, 4).Value = 17600
$ws.Cells.Item(89, 5).Value = 28300
$ws.Cells.Item(89, 6).Value = 22800
$ws.Cells.Item(89, 7).Value = 22600
$ws.Cells.Item(89, 8).Value = 13100
$ws.Cells.Item(89, 9).Value = 31000
$ws.Cells.Item(89, 10).Value = 28700
$ws.Cells.Item(89, 11).Value = 28200
$ws.Cells.Item(89, 12).Value = 37600
$ws.Cells.Item(89, 13).Value = 21600
$ws.Cells.Item(91, 4).Value = -2500
$ws.Cells.Item(91, 5).Value = -2800
$ws.Cells.Item(91, 6).Value = -1600
$ws.Cells.Item(91, 7).Value = -1500
$ws.Cells.Item(91, 8).Value = -5200
$ws.Cells.Item(91, 9).Value = 500
$ws.Cells.Item(91, 10).Value = -1500
$ws.Cells.Item(91, 11).Value = -29400
$ws.Cells.Item(91, 12).Value = -67300
$ws.Cells.Item(91, 13).Value = -84900
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = 0
$ws.Cells.Item(94, 4).Value = -40900
$ws.Cells.Item(94, 5).Value = -38600
$ws.Cells.Item(94, 6).Value = -35800
$ws.Cells.Item(94, 7).Value = -31600
$ws.Cells.Item(94, 8).Value = -18400
$ws.Cells.Item(94, 9).Value = -42500
$ws.Cells.Item(94, 10).Value = -22000
$ws.Cells.Item(94, 11).Value = -31000
$ws.Cells.Item(94, 12).Value = -35500
$ws.Cells.Item(94, 13).Value = -22900
$ws.Cells.Item(96, 4).Value = -4500
$ws.Cells.Item(96, 5).Value = -5800
$ws.Cells.Item(96, 6).Value = -5800
$ws.Cells.Item(96, 7).Value = -5800
$ws.Cells.Item(96, 8).Value = -8000
$ws.Cells.Item(96, 9).Value = -4500
$ws.Cells.Item(96, 10).Value = -4500
$ws.Cells.Item(96, 11).Value = -4500
$ws.Cells.Item(96, 12).Value = -4100
$ws.Cells.Item(96, 13).Value = -4100
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 0
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = 0
$ws.Cells.Item(100, 4).Value = 430600
$ws.Cells.Item(100, 5).Value = 14700
$ws.Cells.Item(100, 6).Value = 14900
$ws.Cells.Item(100, 7).Value = 8100
$ws.Cells.Item(100, 8).Value = 5500
$ws.Cells.Item(100, 9).Value = 9800
$ws.Cells.Item(100, 10).Value = -4600
$ws.Cells.Item(100, 11).Value = -15500
$ws.Cells.Item(100, 12).Value = 5500
$ws.Cells.Item(100, 13).Value = 6500
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(101, 5).Value = 0
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 0
$ws.Cells.Item(101, 9).Value = 0
$ws.Cells.Item(101, 10).Value = 0
$ws.Cells.Item(101, 11).Value = 0
$ws.Cells.Item(101, 12).Value = 0
$ws.Cells.Item(101, 13).Value = 0
$ws.Cells.Item(102, 4).Value = 407400
$ws.Cells.Item(102, 5).Value = 4400
$ws.Cells.Item(102, 6).Value = 2000
$ws.Cells.Item(102, 7).Value = -800
$ws.Cells.Item(102, 8).Value = 200
$ws.Cells.Item(102, 9).Value = -1700
$ws.Cells.Item(102, 10).Value = 2100
$ws.Cells.Item(102, 11).Value = -18300
$ws.Cells.Item(102, 12).Value = 7700
$ws.Cells.Item(102, 13).Value = 5100
